$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 167, shifting existing rows 167-204 down to 168-205
$ws.Rows.Item(167).Insert()

# Populate the newly inserted row 167 with the new data record
$ws.Cells.Item(167, 1).Value = 10
$ws.Cells.Item(167, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(167, 3).Value = "La Araucanía"
$ws.Cells.Item(167, 4).Value = 44855
$ws.Cells.Item(167, 5).Value = 9
$ws.Cells.Item(167, 6).Value = 100114007
$ws.Cells.Item(167, 7).Value = "Jengibre"
$ws.Cells.Item(167, 8).Value = "Sin especificar"
$ws.Cells.Item(167, 9).Value = "Primera"
$ws.Cells.Item(167, 10).Value = 30
$ws.Cells.Item(167, 11).Value = 20000
$ws.Cells.Item(167, 12).Value = 20000
$ws.Cells.Item(167, 13).Value = 20000
$ws.Cells.Item(167, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(167, 15).Value = "Perú"
$ws.Cells.Item(167, 16).Value = 1538
$ws.Cells.Item(167, 17).Value = 13
$ws.Cells.Item(167, 18).Value = "Hortaliza"
